# Applies:
#  1. Table style change on the table (graphicFrame) on slide 16
#     from {E67C4561-71B2-47AB-8AA2-BF868D6A6F11} to {E8E6F3AE-6A1F-4A90-BFCF-72413FED7A14}.
#  2. Re-colour the deck's theme colour scheme from the "Integral" palette
#     to the "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{E8E6F3AE-6A1F-4A90-BFCF-72413FED7A14}")
    }
}

# --- 2. Theme colours -------------------------------------------------------
# Office Theme colour scheme (hex -> decimal BGR-packed RGB used by the OM).
$officeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388     # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501      # accent2  ED7D31
    7  = 10855845     # accent3  A5A5A5
    8  = 49407        # accent4  FFC000
    9  = 12874308     # accent5  4472C4
    10 = 4697456      # accent6  70AD47
    11 = 12673797     # hlink    0563C1
    12 = 7491477      # folHlink 954F72
}

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColorScheme.Colors($i).RGB = $officeColors[$i]
}
